$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.90601666960259
$ws.Range("C2").Value = 9.563708855824501
$ws.Range("D2").Value = 6.641770795019138
$ws.Range("F2").Value = 33.83845325315806
$ws.Range("G2").Value = 45.50837049108168
$ws.Range("H2").Value = 18.27211043341456
$ws.Range("K2").Value = 9.85959058875328
$ws.Range("L2").Value = 10.82934092272948
$ws.Range("B3").Value = 13.70895798143026
$ws.Range("C3").Value = 9.562569478288975
$ws.Range("D3").Value = 6.638967512032064
$ws.Range("F3").Value = 33.72730132999728
$ws.Range("G3").Value = 45.28612249149185
$ws.Range("H3").Value = 18.28076485178492
$ws.Range("K3").Value = 9.735651212971568
$ws.Range("L3").Value = 10.80667904982892
$ws.Range("B4").Value = 13.59049906277789
$ws.Range("C4").Value = 9.561814289662362
$ws.Range("D4").Value = 6.637068113447858
$ws.Range("F4").Value = 33.66628472350207
$ws.Range("G4").Value = 45.15989683925439
$ws.Range("H4").Value = 18.28894047597919
$ws.Range("K4").Value = 9.661492602615105
$ws.Range("L4").Value = 10.79499181844893
$ws.Range("B5").Value = 13.54292857096887
$ws.Range("C5").Value = 9.5614925804722
$ws.Range("D5").Value = 6.636249367309471
$ws.Range("F5").Value = 33.64325221764523
$ws.Range("G5").Value = 45.11106573494369
$ws.Range("H5").Value = 18.29299130375599
$ws.Range("K5").Value = 9.631799136076832
$ws.Range("L5").Value = 10.79079240422311
$ws.Range("B6").Value = 13.53507382478496
$ws.Range("C6").Value = 9.561438318119007
$ws.Range("D6").Value = 6.636110716102428
$ws.Range("F6").Value = 33.63953871680248
$ws.Range("G6").Value = 45.10311570495834
$ws.Range("H6").Value = 18.29370736161548
$ws.Range("K6").Value = 9.626901487334223
$ws.Range("L6").Value = 10.79012919517788
$ws.Range("B7").Value = 13.58985458146373
$ws.Range("C7").Value = 9.561810007491905
$ws.Range("D7").Value = 6.6370572525559
$ws.Range("F7").Value = 33.66596666358942
$ws.Range("G7").Value = 45.15922769012164
$ws.Range("H7").Value = 18.28899219563061
$ws.Range("K7").Value = 9.661089962799339
$ws.Range("L7").Value = 10.79493289937301
$ws.Range("B8").Value = 13.83758371822487
$ws.Range("C8").Value = 9.563327550331737
$ws.Range("D8").Value = 6.640841247872626
$ws.Range("F8").Value = 33.79863623545864
$ws.Range("G8").Value = 45.4296360127069
$ws.Range("H8").Value = 18.27450031961035
$ws.Range("K8").Value = 9.816477056164942
$ws.Range("L8").Value = 10.82106688703852
$ws.Range("B9").Value = 14.34060974372662
$ws.Range("C9").Value = 9.565861489405018
$ws.Range("D9").Value = 6.646844594124901
$ws.Range("F9").Value = 34.11547212460589
$ws.Range("G9").Value = 46.03949028960817
$ws.Range("H9").Value = 18.26880111034103
$ws.Range("K9").Value = 10.13484445608581
$ws.Range("L9").Value = 10.88982799609521
$ws.Range("B10").Value = 14.71679611358918
$ws.Range("C10").Value = 9.567452863660282
$ws.Range("D10").Value = 6.650386913190556
$ws.Range("F10").Value = 34.38172584732997
$ws.Range("G10").Value = 46.53364155017353
$ws.Range("H10").Value = 18.27847160212501
$ws.Range("K10").Value = 10.37470316794128
$ws.Range("L10").Value = 10.95077914278526
$ws.Range("B11").Value = 14.88855313397358
$ws.Range("C11").Value = 9.568117971515289
$ws.Range("D11").Value = 6.651808945624849
$ws.Range("F11").Value = 34.50985038144575
$ws.Range("G11").Value = 46.76786484084246
$ws.Range("H11").Value = 18.28587700969348
$ws.Range("K11").Value = 10.48461028881231
$ws.Range("L11").Value = 10.98071092795414
$ws.Range("B12").Value = 14.95361693779715
$ws.Range("C12").Value = 9.568361361313785
$ws.Range("D12").Value = 6.652320127576251
$ws.Range("F12").Value = 34.55934924851679
$ws.Range("G12").Value = 46.85786112468602
$ws.Range("H12").Value = 18.28911279264847
$ws.Range("K12").Value = 10.52630196869925
$ws.Range("L12").Value = 10.99235638683939
$ws.Range("B13").Value = 14.93960430313955
$ws.Range("C13").Value = 9.568309319976699
$ws.Range("D13").Value = 6.652211251254056
$ws.Range("F13").Value = 34.54864559422903
$ws.Range("G13").Value = 46.83842185922886
$ws.Range("H13").Value = 18.28839673219775
$ws.Range("K13").Value = 10.51732036584138
$ws.Range("L13").Value = 10.98983460158447
$ws.Range("B14").Value = 14.89390588632661
$ws.Range("C14").Value = 9.568138165210648
$ws.Range("D14").Value = 6.651851548656946
$ws.Range("F14").Value = 34.51390318157036
$ws.Range("G14").Value = 46.77524314082586
$ws.Range("H14").Value = 18.28613457664965
$ws.Range("K14").Value = 10.48803907803586
$ws.Range("L14").Value = 10.98166281166905
$ws.Range("B15").Value = 14.86591533980929
$ws.Range("C15").Value = 9.568032224257541
$ws.Range("D15").Value = 6.65162766149556
$ws.Range("F15").Value = 34.49274937995921
$ws.Range("G15").Value = 46.73671209698684
$ws.Range("H15").Value = 18.28480510970791
$ws.Range("K15").Value = 10.47011162000705
$ws.Range("L15").Value = 10.97669765942228
$ws.Range("B16").Value = 14.705578242918
$ws.Range("C16").Value = 9.567408212494598
$ws.Range("D16").Value = 6.650290162415811
$ws.Range("F16").Value = 34.37349134246157
$ws.Range("G16").Value = 46.51851963273438
$ws.Range("H16").Value = 18.27804807391743
$ws.Range("K16").Value = 10.36753285487621
$ws.Range("L16").Value = 10.94886688459751
$ws.Range("B17").Value = 14.60732912446157
$ws.Range("C17").Value = 9.567010326443045
$ws.Range("D17").Value = 6.649421087815398
$ws.Range("F17").Value = 34.30210621116475
$ws.Range("G17").Value = 46.38704281355071
$ws.Range("H17").Value = 18.27467235000528
$ws.Range("K17").Value = 10.30477741587892
$ws.Range("L17").Value = 10.93235410500508
$ws.Range("B18").Value = 14.55088185276389
$ws.Range("C18").Value = 9.566775937596427
$ws.Range("D18").Value = 6.648903397559754
$ws.Range("F18").Value = 34.26170834783186
$ws.Range("G18").Value = 46.3123120391528
$ws.Range("H18").Value = 18.27301374804433
$ws.Range("K18").Value = 10.26875936951054
$ws.Range("L18").Value = 10.92306431161944
$ws.Range("B19").Value = 14.53178267596033
$ws.Range("C19").Value = 9.566695627330937
$ws.Range("D19").Value = 6.648725058085045
$ws.Range("F19").Value = 34.24814462246684
$ws.Range("G19").Value = 46.28716418796017
$ws.Range("H19").Value = 18.27250080343671
$ws.Range("K19").Value = 10.25657883229374
$ws.Range("L19").Value = 10.91995483657203
$ws.Range("B20").Value = 14.61778188044557
$ws.Range("C20").Value = 9.567053255085657
$ws.Range("D20").Value = 6.649515447306935
$ws.Range("F20").Value = 34.30963705687006
$ws.Range("G20").Value = 46.40094689940265
$ws.Range("H20").Value = 18.27500241434391
$ws.Range("K20").Value = 10.31145015379808
$ws.Range("L20").Value = 10.93409044240739
$ws.Range("B21").Value = 14.90732851673239
$ws.Range("C21").Value = 9.568188667551624
$ws.Range("D21").Value = 6.651957943840508
$ws.Range("F21").Value = 34.52408146776609
$ws.Range("G21").Value = 46.79376540607552
$ws.Range("H21").Value = 18.28678732281528
$ws.Range("K21").Value = 10.49663807180932
$ws.Range("L21").Value = 10.98405467643009
$ws.Range("B22").Value = 15.09666915468084
$ws.Range("C22").Value = 9.568881337798645
$ws.Range("D22").Value = 6.653395033847621
$ws.Range("F22").Value = 34.66993723999553
$ws.Range("G22").Value = 47.05805086384962
$ws.Range("H22").Value = 18.2970041062758
$ws.Range("K22").Value = 10.61807131797877
$ws.Range("L22").Value = 11.01851845949242
$ws.Range("B23").Value = 14.99562655305525
$ws.Range("C23").Value = 9.568516169924743
$ws.Range("D23").Value = 6.652642625646513
$ws.Range("F23").Value = 34.5915784966024
$ws.Range("G23").Value = 46.91632428884689
$ws.Range("H23").Value = 18.29132144754739
$ws.Range("K23").Value = 10.55323694095004
$ws.Range("L23").Value = 10.99996107117756
$ws.Range("B24").Value = 14.61305606764371
$ws.Range("C24").Value = 9.567033864609968
$ws.Range("D24").Value = 6.64947284356478
$ws.Range("F24").Value = 34.3062303584098
$ws.Range("G24").Value = 46.39465818839529
$ws.Range("H24").Value = 18.27485231326141
$ws.Range("K24").Value = 10.3084332160941
$ws.Range("L24").Value = 10.93330480955018
$ws.Range("B25").Value = 14.20308610811449
$ws.Range("C25").Value = 9.565223403902193
$ws.Range("D25").Value = 6.645372707675763
$ws.Range("F25").Value = 34.02379541384769
$ws.Range("G25").Value = 45.86622496038454
$ws.Range("H25").Value = 18.26790893713802
$ws.Range("K25").Value = 10.04749997677295
$ws.Range("L25").Value = 10.86937470185336
